$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting rows 8-13 down to 9-14
$ws.Rows.Item(8).Insert()

# Fill the new row 8 with TIME_STEP variable info
$ws.Cells.Item(8, 1).Value = "TIME_STEP"
$ws.Cells.Item(8, 2).Value = "s"
$ws.Cells.Item(8, 3).Value = "float"
$ws.Cells.Item(8, 4).Value = "User defined time step for the thermal-hydraulic loop. Should be in the range [STPMIN, STPMAX] (boundary included). "
$ws.Cells.Item(8, 5).Value = 0.3

# Match row height of sibling rows (STPMIN/STPMAX/etc use ht=29)
$ws.Rows.Item(8).RowHeight = 29

# Restore freeze panes (xSplit=1, ySplit=2) and update the active selection
$aw = $excel.ActiveWindow
$aw.FreezePanes = $false
$ws.Range("B3").Select()
$aw.FreezePanes = $true
$ws.Range("E9").Select()

Write-Host "Done"
